# feat(master-item-import): Include color info when importing master measures from Excel file
# Implements #245
#
# Adds two new shared-string JSON blobs (single color info, and multi-color
# breakpoint info) and wires them into new cells G2/H2 (row 2), G3 (row 3),
# and H4 (row 4) on the "Sales" sheet, adjusting row heights + the active
# selection to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$singleColor = '{
  "color": "#8a85c6",
  "index": 8
}'

$multiColor = '{
  "colors": [
    {
      "color": "#006580",
      "index": 6
    },
    {
      "color": "#ac4d58",
      "index": 10
    },
    {
      "color": "#4477aa",
      "index": -1
    },
    {
      "color": "#7db8da",
      "index": -1
    }
  ],
  "breakTypes": [
    true,
    true,
    false
  ],
  "limits": [
    0.223,
    0.491,
    0.728
  ],
  "limitType": "percent"
}'

# Row 2: single color (G2) + multi-color (H2)
$ws.Range("G2").Value = $singleColor
$ws.Range("H2").Value = $multiColor
$ws.Rows.Item(2).RowHeight = 188

# Row 3: single color (G3)
$ws.Range("G3").Value = $singleColor
$ws.Rows.Item(3).RowHeight = 60

# Row 4: multi-color (H4)
$ws.Range("H4").Value = $multiColor
$ws.Rows.Item(4).RowHeight = 190

# Update the view: scroll/selection moved from A11/A12 to B2/H3
$ws.Application.ActiveWindow.ScrollRow = 2
$ws.Application.ActiveWindow.ScrollColumn = 2
$ws.Range("H3").Select()
